# Workbook: testes/teste_pre_dot.xlsx
# Commit: "Criadas as funcoes de extracao da capa de dados, de extracao das
#          partes representantes, de processos relacionados. A funcao de
#          screenshot esta 30% concluida."
#
# The single worksheet (Planilha1) holds a work queue of CNPJ numbers in
# column A (A1 = "nada", A2:A12 = CNPJ strings). The row that was sitting in
# A2 has now been processed, so it is removed from the top of the queue and
# appended to the bottom - every other row shifts up by one position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rotate A2:A12: move the old A2 value to the end, shifting the rest up.
$movedValue = $ws.Range("A2").Value()
$ws.Range("A2:A11").Value = $ws.Range("A3:A12").Value()
$ws.Range("A12").Value = $movedValue

# --- Column A grew a bit wider (content lengths now fit differently).
$ws.Columns("A").ColumnWidth = 17

# --- The user's active selection moved to E8 (single cell).
$ws.Range("E8").Select()
